$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the workbook window (matches xWindow 0 -> 1340 in the saved file).
$wb.Windows.Item(1).Left = 1340

# Write new strings in the order they should first appear in sharedStrings.xml:
# Whale Age (J1), 20+ (J2), 17+ (J9), Gear length info.. (I14), pApt (K1), Y (K2)

$ws.Range("J1").Value = "Whale Age"
$ws.Range("J2").Value = "20+"
$ws.Range("J9").Value = "17+"
$ws.Range("I14").Value = "Gear length info from PCCS Page; Lobster trap"
$ws.Range("K1").Value = "pApt"
$ws.Range("K2").Value = "Y"

# Remaining J column numeric/shared values
$ws.Range("J3").Value = 18
$ws.Range("J4").Value = 8
$ws.Range("J5").Value = 7
$ws.Range("J6").Value = 3
$ws.Range("J7").Value = 8
$ws.Range("J8").Value = 7
$ws.Range("J10").Value = 2
$ws.Range("J11").Value = 1
$ws.Range("J12").Value = 1
$ws.Range("J13").Value = 5
$ws.Range("J14").Value = 1

# Remaining K column (all "Y")
$ws.Range("K3").Value = "Y"
$ws.Range("K4").Value = "Y"
$ws.Range("K5").Value = "Y"
$ws.Range("K6").Value = "Y"
$ws.Range("K7").Value = "Y"
$ws.Range("K8").Value = "Y"
$ws.Range("K9").Value = "Y"
$ws.Range("K10").Value = "Y"
$ws.Range("K11").Value = "Y"
$ws.Range("K12").Value = "Y"
$ws.Range("K13").Value = "Y"
$ws.Range("K14").Value = "Y"

$ws.Range("K15").Select()
